$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 329
$ws.Range("F4").Value = 637
$ws.Range("F5").Value = 27
$ws.Range("F6").Value = 461
$ws.Range("F8").Value = 2125
$ws.Range("F9").Value = 874
$ws.Range("F10").Value = 834
$ws.Range("F11").Value = 404
$ws.Range("F12").Value = 79
$ws.Range("F13").Value = 427
$ws.Range("F14").Value = 321
$ws.Range("F16").Value = 887
$ws.Range("F18").Value = 32
$ws.Range("F19").Value = 1717
$ws.Range("F21").Value = 29
$ws.Range("F23").Value = 59
$ws.Range("F24").Value = 506
$ws.Range("F25").Value = 1454
$ws.Range("F27").Value = 523
$ws.Range("F29").Value = 593
$ws.Range("F30").Value = 415
$ws.Range("F31").Value = 2376
$ws.Range("F32").Value = 380
$ws.Range("F33").Value = 90
$ws.Range("F35").Value = 598
$ws.Range("F36").Value = 474
$ws.Range("F37").Value = 183
$ws.Range("F38").Value = 916
$ws.Range("F39").Value = 696
$ws.Range("F41").Value = 449
$ws.Range("F42").Value = 415

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 79
$ws.Range("F11").Value = 52
$ws.Range("F23").Value = 101
$ws.Range("F24").Value = 434

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2919
$ws.Range("F6").Value = 311

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 637
$ws.Range("F8").Value = 27
$ws.Range("F10").Value = 461
$ws.Range("F11").Value = 874
$ws.Range("F12").Value = 834
$ws.Range("F13").Value = 404
$ws.Range("F14").Value = 79
$ws.Range("F15").Value = 321
$ws.Range("F18").Value = 887
$ws.Range("F19").Value = 79
$ws.Range("F21").Value = 32
$ws.Range("F22").Value = 311
$ws.Range("F23").Value = 1719
$ws.Range("F25").Value = 29
$ws.Range("F27").Value = 52
$ws.Range("F29").Value = 506
$ws.Range("F31").Value = 1454
$ws.Range("F34").Value = 523
$ws.Range("F36").Value = 593
$ws.Range("F37").Value = 415
$ws.Range("F38").Value = 90
$ws.Range("F40").Value = 474
$ws.Range("F41").Value = 183
$ws.Range("F42").Value = 916
$ws.Range("F44").Value = 101
$ws.Range("F45").Value = 434
$ws.Range("F46").Value = 696
$ws.Range("F48").Value = 449
$ws.Range("F49").Value = 416

